$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "End of class Feb-03": remove Matt Klint (Clark State College) from the
# roster. He occupied the two-row record at rows 15:16 (name row + email/
# school row). Deleting the rows shifts everyone below up by one record.
$ws.Range("A15:A16").EntireRow.Delete()

# Reflect the print area shrinking along with the roster (30 -> 28 rows).
$ws.PageSetup.PrintArea = '$A$1:$M$28'

# Leave the selection where the instructor made the edit.
$ws.Range("A15:H16").Select()
